# Applies the "Updated symbol list" data refresh for cryptos.xlsx (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin name / link swaps in rows 9-14 (plain text columns, no numeric coercion risk) ---
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'

# --- "Hora" column: every data row moves from 15 to 16 (one range-wide text write) ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "16"

# --- Price (D) / Volume(1h) (E) refreshed numbers, kept as Text so "0.05901"/"5.47%" etc. survive verbatim ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '258.57'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '5.47%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.25'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-1.98%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.222'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-0.41%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05901'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.29%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.711'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.49%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8661'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.51%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.002'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '15.91%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1411'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.49%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07171'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.87%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.03161'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.02%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09226'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.71%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001548'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.32%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006053'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-93.98%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005880'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.68%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.501'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.32%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.225'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.52%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.204'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.38%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3175'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.33%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03528'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '5.96%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1289'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.95%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.561'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.99%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04190'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.69%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-1.07%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001225'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.09%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004557'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.98%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.13%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001471'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.49%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03839'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.96%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005503'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '6.56%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1101'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.42%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002372'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-3.14%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01083'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '15.98%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005411'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.19%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.13%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.09489'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '26.51%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002135'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-12.58%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.13%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.13%'
